$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Has2"
$ws.Range("C2").Value = "Hmmr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.772687
$ws.Range("H2").Value = 2.318061
$ws.Range("I2").Value = 0.02743927362033742
$ws.Range("J2").Value = 0.02743927362033742
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.26464
$ws.Range("N2").Value = 6.79392
$ws.Range("O2").Value = 0.1949446419849994
$ws.Range("P2").Value = 0.1949446419849994
$ws.Range("Q2").Value = 1.74985788768
$ws.Range("R2").Value = 15.74872098912
$ws.Range("S2").Value = 0.005349139372245117
$ws.Range("T2").Value = 0.005349139372245117

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Has2"
$ws.Range("C3").Value = "Hmmr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.772687
$ws.Range("H3").Value = 2.318061
$ws.Range("I3").Value = 0.02743927362033742
$ws.Range("J3").Value = 0.02743927362033742
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.432364333333333
$ws.Range("N3").Value = 4.297093
$ws.Range("O3").Value = 0.1233007242448023
$ws.Range("P3").Value = 0.1233007242448023
$ws.Range("Q3").Value = 1.106769299630333
$ws.Range("R3").Value = 9.960923696673001
$ws.Range("S3").Value = 0.003383282310138901
$ws.Range("T3").Value = 0.003383282310138902

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Has2"
$ws.Range("C4").Value = "Hmmr"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.772687
$ws.Range("H4").Value = 2.318061
$ws.Range("I4").Value = 0.02743927362033742
$ws.Range("J4").Value = 0.02743927362033742
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.324070333333334
$ws.Range("N4").Value = 6.972211000000001
$ws.Range("O4").Value = 0.2000605213542218
$ws.Range("P4").Value = 0.2000605213542218
$ws.Range("Q4").Value = 1.795778933652334
$ws.Range("R4").Value = 16.162010402871
$ws.Range("S4").Value = 0.005489515386065851
$ws.Range("T4").Value = 0.005489515386065852

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Has2"
$ws.Range("C5").Value = "Hmmr"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.772687
$ws.Range("H5").Value = 2.318061
$ws.Range("I5").Value = 0.02743927362033742
$ws.Range("J5").Value = 0.02743927362033742
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.595761666666667
$ws.Range("N5").Value = 16.787285
$ws.Range("O5").Value = 0.4816941124159764
$ws.Range("P5").Value = 0.4816941124159765
$ws.Range("Q5").Value = 4.323772294931667
$ws.Range("R5").Value = 38.913950654385
$ws.Range("S5").Value = 0.01321733655188755
$ws.Range("T5").Value = 0.01321733655188755

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Has2"
$ws.Range("C6").Value = "Hmmr"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 24.47381466666667
$ws.Range("H6").Value = 73.42144400000001
$ws.Range("I6").Value = 0.86910184482474
$ws.Range("J6").Value = 0.86910184482474
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 2.26464
$ws.Range("N6").Value = 6.79392
$ws.Range("O6").Value = 0.1949446419849994
$ws.Range("P6").Value = 0.1949446419849994
$ws.Range("Q6").Value = 55.42437964672001
$ws.Range("R6").Value = 498.8194168204801
$ws.Range("S6").Value = 0.1694267479878614
$ws.Range("T6").Value = 0.1694267479878614

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Has2"
$ws.Range("C7").Value = "Hmmr"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 24.47381466666667
$ws.Range("H7").Value = 73.42144400000001
$ws.Range("I7").Value = 0.86910184482474
$ws.Range("J7").Value = 0.86910184482474
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.432364333333333
$ws.Range("N7").Value = 4.297093
$ws.Range("O7").Value = 0.1233007242448023
$ws.Range("P7").Value = 0.1233007242448023
$ws.Range("Q7").Value = 35.05541922914356
$ws.Range("R7").Value = 315.4987730622921
$ws.Range("S7").Value = 0.1071608869093842
$ws.Range("T7").Value = 0.1071608869093842

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Has2"
$ws.Range("C8").Value = "Hmmr"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 24.47381466666667
$ws.Range("H8").Value = 73.42144400000001
$ws.Range("I8").Value = 0.86910184482474
$ws.Range("J8").Value = 0.86910184482474
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.324070333333334
$ws.Range("N8").Value = 6.972211000000001
$ws.Range("O8").Value = 0.2000605213542218
$ws.Range("P8").Value = 0.2000605213542218
$ws.Range("Q8").Value = 56.87886661029824
$ws.Range("R8").Value = 511.9097994926841
$ws.Range("S8").Value = 0.1738729681855535
$ws.Range("T8").Value = 0.1738729681855535

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Has2"
$ws.Range("C9").Value = "Hmmr"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 24.47381466666667
$ws.Range("H9").Value = 73.42144400000001
$ws.Range("I9").Value = 0.86910184482474
$ws.Range("J9").Value = 0.86910184482474
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.595761666666667
$ws.Range("N9").Value = 16.787285
$ws.Range("O9").Value = 0.4816941124159764
$ws.Range("P9").Value = 0.4816941124159765
$ws.Range("Q9").Value = 136.9496339488378
$ws.Range("R9").Value = 1232.54670553954
$ws.Range("S9").Value = 0.4186412417419408
$ws.Range("T9").Value = 0.4186412417419408

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Has2"
$ws.Range("C10").Value = "Hmmr"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.913391
$ws.Range("H10").Value = 8.740173
$ws.Range("I10").Value = 0.1034588815549226
$ws.Range("J10").Value = 0.1034588815549226
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 2.26464
$ws.Range("N10").Value = 6.79392
$ws.Range("O10").Value = 0.1949446419849994
$ws.Range("P10").Value = 0.1949446419849994
$ws.Range("Q10").Value = 6.59778179424
$ws.Range("R10").Value = 59.38003614816
$ws.Range("S10").Value = 0.02016875462489284
$ws.Range("T10").Value = 0.02016875462489284

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Has2"
$ws.Range("C11").Value = "Hmmr"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.913391
$ws.Range("H11").Value = 8.740173
$ws.Range("I11").Value = 0.1034588815549226
$ws.Range("J11").Value = 0.1034588815549226
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.432364333333333
$ws.Range("N11").Value = 4.297093
$ws.Range("O11").Value = 0.1233007242448023
$ws.Range("P11").Value = 0.1233007242448023
$ws.Range("Q11").Value = 4.173037357454334
$ws.Range("R11").Value = 37.557336217089
$ws.Range("S11").Value = 0.01275655502527917
$ws.Range("T11").Value = 0.01275655502527917

$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Has2"
$ws.Range("C12").Value = "Hmmr"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.913391
$ws.Range("H12").Value = 8.740173
$ws.Range("I12").Value = 0.1034588815549226
$ws.Range("J12").Value = 0.1034588815549226
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.324070333333334
$ws.Range("N12").Value = 6.972211000000001
$ws.Range("O12").Value = 0.2000605213542218
$ws.Range("P12").Value = 0.2000605213542218
$ws.Range("Q12").Value = 6.770925592500335
$ws.Range("R12").Value = 60.93833033250301
$ws.Range("S12").Value = 0.0206980377826025
$ws.Range("T12").Value = 0.0206980377826025

$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Has2"
$ws.Range("C13").Value = "Hmmr"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.913391
$ws.Range("H13").Value = 8.740173
$ws.Range("I13").Value = 0.1034588815549226
$ws.Range("J13").Value = 0.1034588815549226
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 5.595761666666667
$ws.Range("N13").Value = 16.787285
$ws.Range("O13").Value = 0.4816941124159764
$ws.Range("P13").Value = 0.4816941124159765
$ws.Range("Q13").Value = 16.30264167781167
$ws.Range("R13").Value = 146.723775100305
$ws.Range("S13").Value = 0.04983553412214807
$ws.Range("T13").Value = 0.04983553412214808
